$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1232.5883
$ws.Range("I33").Value = 1163.4667
$ws.Range("K33").Value = 1163.4667
$ws.Range("M33").Value = -934.4666999999999
$ws.Range("H69").Value = 8333.333000000001
$ws.Range("I69").Value = 9000
$ws.Range("J69").Value = 7000
$ws.Range("K69").Value = 27000
$ws.Range("L69").Value = 21000
$ws.Range("M69").Value = -26126
$ws.Range("N69").Value = -22748
$ws.Range("H72").Value = 8333.333000000001
$ws.Range("I72").Value = 9000
$ws.Range("J72").Value = 7000
$ws.Range("K72").Value = 81000
$ws.Range("L72").Value = 63000
$ws.Range("M72").Value = -76632
$ws.Range("N72").Value = -71736
$ws.Range("H100").Value = 3491.7058
$ws.Range("I100").Value = 3050.923
$ws.Range("K100").Value = 3050.923
$ws.Range("M100").Value = -2509.923
$ws.Range("H107").Value = 730.619
$ws.Range("I107").Value = 804.8889
$ws.Range("K107").Value = 804.8889
$ws.Range("M107").Value = 1115.1111
$ws.Range("H113").Value = 4721.4736
$ws.Range("I113").Value = 4101.25
$ws.Range("J113").Value = 5172.5454
$ws.Range("K113").Value = 4101.25
$ws.Range("L113").Value = 5172.5454
$ws.Range("M113").Value = -847.25
$ws.Range("N113").Value = -11680.5454
$ws.Range("H132").Value = 1293.6595
$ws.Range("I132").Value = 881.4651
$ws.Range("K132").Value = 2644.3953
$ws.Range("M132").Value = -114.3953000000001
$ws.Range("H137").Value = 2442.5386
$ws.Range("I137").Value = 2101.1482
$ws.Range("J137").Value = 3210.6667
$ws.Range("K137").Value = 6303.444600000001
$ws.Range("L137").Value = 9632.000100000001
$ws.Range("M137").Value = -3753.444600000001
$ws.Range("N137").Value = -14732.0001
$ws.Range("H138").Value = 2531.9697
$ws.Range("J138").Value = 3309.3684
$ws.Range("L138").Value = 9928.1052
$ws.Range("N138").Value = -20208.1052

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 201.94737
$ws.Range("I5").Value = 167.63637
$ws.Range("K5").Value = 167.63637
$ws.Range("M5").Value = -55.63637
$ws.Range("H32").Value = 64422
$ws.Range("I32").Value = 45738.082
$ws.Range("K32").Value = 45738.082
$ws.Range("M32").Value = -45451.082
$ws.Range("H41").Value = 9961
$ws.Range("I41").Value = 3281.3333
$ws.Range("K41").Value = 3281.3333
$ws.Range("M41").Value = -2867.3333
$ws.Range("H45").Value = 391468
$ws.Range("I45").Value = 596755.4
$ws.Range("K45").Value = 596755.4
$ws.Range("M45").Value = -596378.4
$ws.Range("H132").Value = 2036.7142
$ws.Range("I132").Value = 1116.6
$ws.Range("K132").Value = 3349.8
$ws.Range("M132").Value = -819.7999999999997

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 201.94737
$ws.Range("I4").Value = 167.63637
$ws.Range("K4").Value = 167.63637
$ws.Range("M4").Value = -52.63637
$ws.Range("H20").Value = 10683
$ws.Range("I20").Value = 9244.066000000001
$ws.Range("J20").Value = 12645.182
$ws.Range("K20").Value = 9244.066000000001
$ws.Range("L20").Value = 12645.182
$ws.Range("M20").Value = -8997.066000000001
$ws.Range("N20").Value = -13139.182
$ws.Range("H99").Value = 2479.25
$ws.Range("I99").Value = 2047
$ws.Range("J99").Value = 3199.6667
$ws.Range("K99").Value = 2047
$ws.Range("L99").Value = 3199.6667
$ws.Range("M99").Value = -549
$ws.Range("N99").Value = -6195.6667
$ws.Range("H107").Value = 1364.75
$ws.Range("I107").Value = 1416.8572
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1416.8572
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 503.1428000000001
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 1411.8485
$ws.Range("I134").Value = 1380.3549
$ws.Range("K134").Value = 4141.0647
$ws.Range("M134").Value = -1606.0647

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1128.64
$ws.Range("I107").Value = 821.8570999999999
$ws.Range("J107").Value = 1519.091
$ws.Range("K107").Value = 821.8570999999999
$ws.Range("L107").Value = 1519.091
$ws.Range("M107").Value = 1098.1429
$ws.Range("N107").Value = -5359.091
$ws.Range("H122").Value = 3276.2778
$ws.Range("I122").Value = 2229.8
$ws.Range("J122").Value = 4584.375
$ws.Range("K122").Value = 6689.400000000001
$ws.Range("L122").Value = 13753.125
$ws.Range("M122").Value = -4239.400000000001
$ws.Range("N122").Value = -18653.125
$ws.Range("H134").Value = 1854.871
$ws.Range("I134").Value = 1783.3667
$ws.Range("K134").Value = 5350.1001
$ws.Range("M134").Value = -2815.1001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 9.232324
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 9.204082
$ws.Range("K2").Value = 72
$ws.Range("L2").Value = 55.224492
$ws.Range("M2").Value = 41
$ws.Range("N2").Value = -281.224492
$ws.Range("H11").Value = 1668166
$ws.Range("I11").Value = 999
$ws.Range("J11").Value = 2501749.5
$ws.Range("K11").Value = 2997
$ws.Range("L11").Value = 7505248.5
$ws.Range("M11").Value = -2857
$ws.Range("N11").Value = -7505528.5
$ws.Range("H26").Value = 339.81818
$ws.Range("I26").Value = 314.14285
$ws.Range("J26").Value = 384.75
$ws.Range("K26").Value = 942.4285500000001
$ws.Range("L26").Value = 1154.25
$ws.Range("M26").Value = -654.4285500000001
$ws.Range("N26").Value = -1730.25
$ws.Range("H131").Value = 24648.076
$ws.Range("I131").Value = 939.8
$ws.Range("J131").Value = 39465.75
$ws.Range("K131").Value = 2819.4
$ws.Range("L131").Value = 118397.25
$ws.Range("M131").Value = 2220.6
$ws.Range("N131").Value = -128477.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2416.2942
$ws.Range("I102").Value = 2429.8125
$ws.Range("J102").Value = 2200
$ws.Range("K102").Value = 2429.8125
$ws.Range("L102").Value = 2200
$ws.Range("M102").Value = -807.8125
$ws.Range("N102").Value = -5444
$ws.Range("H113").Value = 2494.4167
$ws.Range("I113").Value = 2605.4285
$ws.Range("J113").Value = 2339
$ws.Range("K113").Value = 2605.4285
$ws.Range("L113").Value = 2339
$ws.Range("M113").Value = -435.4285
$ws.Range("N113").Value = -6679

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 37061.125
$ws.Range("J7").Value = 6749.75
$ws.Range("L7").Value = 6749.75
$ws.Range("N7").Value = -6973.75
$ws.Range("H100").Value = 7304.45
$ws.Range("I100").Value = 7396.1875
$ws.Range("K100").Value = 7396.1875
$ws.Range("M100").Value = -6855.1875
$ws.Range("H126").Value = 37061.125
$ws.Range("J126").Value = 6749.75
$ws.Range("L126").Value = 20249.25
$ws.Range("N126").Value = -25189.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 55556460
$ws.Range("I107").Value = 1015.375
$ws.Range("K107").Value = 3046.125
$ws.Range("M107").Value = -1126.125
$ws.Range("H122").Value = 1400.5883
$ws.Range("J122").Value = 1399.1666
$ws.Range("L122").Value = 4197.4998
$ws.Range("N122").Value = -9097.4998
